$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A119").Value = 41960
$ws.Range("B119").Value = 0.54722222222222217
$ws.Range("C119").Value = 0.56388888888888888
$ws.Range("D119").Value = 5
$ws.Range("E119").Formula = "=IF(AND(NOT(ISBLANK(B119)),NOT(ISBLANK(C119))), (C119-B119) * 24 - D119/60, """")"
$ws.Range("F119").Value = "Coding"

$ws.Range("A120").Select()
